# Auto-generated Excel COM-interop edit script
# Applies the diff changes to Stocks (sheet1) and SupportHistory (sheet2)

$wb = $excel.ActiveWorkbook
$wsStocks = $wb.Worksheets.Item("Stocks")
$wsSupport = $wb.Worksheets.Item("SupportHistory")

#### Stocks sheet (sheet1) changes ####
# Row 2
$wsStocks.Cells.Item(2,8).Value = $true
$wsStocks.Cells.Item(2,9).Value = 2273.5
$wsStocks.Cells.Item(2,10).NumberFormat = "@"
$wsStocks.Cells.Item(2,10).Value = "2025-07-14"
$wsStocks.Cells.Item(2,11).Value = 1932.47

# Row 6
$wsStocks.Cells.Item(6,8).Value = $true
$wsStocks.Cells.Item(6,9).Value = 1623.8
$wsStocks.Cells.Item(6,10).NumberFormat = "@"
$wsStocks.Cells.Item(6,10).Value = "2025-07-14"
$wsStocks.Cells.Item(6,11).Value = 1380.23

# Row 7
$wsStocks.Cells.Item(7,8).Value = $true
$wsStocks.Cells.Item(7,9).Value = 2045
$wsStocks.Cells.Item(7,10).NumberFormat = "@"
$wsStocks.Cells.Item(7,10).Value = "2025-07-14"
$wsStocks.Cells.Item(7,11).Value = 1675.25

# Row 10
$wsStocks.Cells.Item(10,8).Value = $true
$wsStocks.Cells.Item(10,9).Value = 1642
$wsStocks.Cells.Item(10,10).NumberFormat = "@"
$wsStocks.Cells.Item(10,10).Value = "2025-07-14"
$wsStocks.Cells.Item(10,11).Value = 1395.7

# Row 11
$wsStocks.Cells.Item(11,10).NumberFormat = "@"
$wsStocks.Cells.Item(11,10).Value = "2025-06-19"
$wsStocks.Cells.Item(11,11).Value = 1140.88

# Row 12
$wsStocks.Cells.Item(12,8).Value = $true
$wsStocks.Cells.Item(12,9).Value = 2277.9
$wsStocks.Cells.Item(12,10).NumberFormat = "@"
$wsStocks.Cells.Item(12,10).Value = "2025-07-14"
$wsStocks.Cells.Item(12,11).Value = 1936.21

# Row 16
$wsStocks.Cells.Item(16,7).Value = $true
$wsStocks.Cells.Item(16,8).Value = $true
$wsStocks.Cells.Item(16,9).Value = 1689
$wsStocks.Cells.Item(16,10).NumberFormat = "@"
$wsStocks.Cells.Item(16,10).Value = "2025-07-14"
$wsStocks.Cells.Item(16,11).Value = 1435.65

# Row 21
$wsStocks.Cells.Item(21,10).NumberFormat = "@"
$wsStocks.Cells.Item(21,10).Value = "2025-06-19"

# Row 22
$wsStocks.Cells.Item(22,8).Value = $true
$wsStocks.Cells.Item(22,9).Value = 1061
$wsStocks.Cells.Item(22,10).NumberFormat = "@"
$wsStocks.Cells.Item(22,10).Value = "2025-07-14"
$wsStocks.Cells.Item(22,11).Value = 901.85

# Row 23
$wsStocks.Cells.Item(23,10).NumberFormat = "@"
$wsStocks.Cells.Item(23,10).Value = "2025-06-19"

# Row 31
$wsStocks.Cells.Item(31,7).Value = $true

# Row 40
$wsStocks.Cells.Item(40,10).NumberFormat = "@"
$wsStocks.Cells.Item(40,10).Value = "2025-06-19"

# Row 45
$wsStocks.Cells.Item(45,10).NumberFormat = "@"
$wsStocks.Cells.Item(45,10).Value = "2025-06-19"

# Row 56
$wsStocks.Cells.Item(56,8).Value = $false
$wsStocks.Cells.Item(56,9).ClearContents()
$wsStocks.Cells.Item(56,10).ClearContents()
$wsStocks.Cells.Item(56,11).ClearContents()

# Row 57
$wsStocks.Cells.Item(57,8).Value = $true
$wsStocks.Cells.Item(57,9).Value = 590
$wsStocks.Cells.Item(57,10).NumberFormat = "@"
$wsStocks.Cells.Item(57,10).Value = "2025-07-14"
$wsStocks.Cells.Item(57,11).Value = 501.5

# Row 75
$wsStocks.Cells.Item(75,8).Value = $true
$wsStocks.Cells.Item(75,9).Value = 590
$wsStocks.Cells.Item(75,10).NumberFormat = "@"
$wsStocks.Cells.Item(75,10).Value = "2025-07-14"
$wsStocks.Cells.Item(75,11).Value = 501.5

# Row 76
$wsStocks.Cells.Item(76,10).NumberFormat = "@"
$wsStocks.Cells.Item(76,10).Value = "2025-06-19"
$wsStocks.Cells.Item(76,11).Value = 430.6

# Row 79
$wsStocks.Cells.Item(79,7).Value = $true
$wsStocks.Cells.Item(79,8).Value = $true
$wsStocks.Cells.Item(79,9).Value = 364
$wsStocks.Cells.Item(79,10).NumberFormat = "@"
$wsStocks.Cells.Item(79,10).Value = "2025-07-14"
$wsStocks.Cells.Item(79,11).Value = 309.4

# Row 81
$wsStocks.Cells.Item(81,8).Value = $false
$wsStocks.Cells.Item(81,9).ClearContents()
$wsStocks.Cells.Item(81,10).ClearContents()
$wsStocks.Cells.Item(81,11).ClearContents()

# Row 82
$wsStocks.Cells.Item(82,9).Value = 408.4
$wsStocks.Cells.Item(82,10).NumberFormat = "@"
$wsStocks.Cells.Item(82,10).Value = "2025-07-14"
$wsStocks.Cells.Item(82,11).Value = 286.14

# Row 84
$wsStocks.Cells.Item(84,10).NumberFormat = "@"
$wsStocks.Cells.Item(84,10).Value = "2025-06-19"
$wsStocks.Cells.Item(84,11).Value = 355.37

# Row 87
$wsStocks.Cells.Item(87,8).Value = $false
$wsStocks.Cells.Item(87,9).ClearContents()
$wsStocks.Cells.Item(87,10).ClearContents()
$wsStocks.Cells.Item(87,11).ClearContents()

# Row 89
$wsStocks.Cells.Item(89,10).NumberFormat = "@"
$wsStocks.Cells.Item(89,10).Value = "2025-06-19"

# Row 90
$wsStocks.Cells.Item(90,10).NumberFormat = "@"
$wsStocks.Cells.Item(90,10).Value = "2025-06-19"
$wsStocks.Cells.Item(90,11).Value = 400.3

# Row 91
$wsStocks.Cells.Item(91,10).NumberFormat = "@"
$wsStocks.Cells.Item(91,10).Value = "2025-06-19"

# Row 92
$wsStocks.Cells.Item(92,9).Value = 333.9
$wsStocks.Cells.Item(92,10).NumberFormat = "@"
$wsStocks.Cells.Item(92,10).Value = "2025-07-14"
$wsStocks.Cells.Item(92,11).Value = 283.81

# Row 94
$wsStocks.Cells.Item(94,10).NumberFormat = "@"
$wsStocks.Cells.Item(94,10).Value = "2025-06-19"

# Row 96
$wsStocks.Cells.Item(96,8).Value = $true
$wsStocks.Cells.Item(96,9).Value = 526.9
$wsStocks.Cells.Item(96,10).NumberFormat = "@"
$wsStocks.Cells.Item(96,10).Value = "2025-07-25"
$wsStocks.Cells.Item(96,11).Value = 443.7

# Row 97
$wsStocks.Cells.Item(97,10).NumberFormat = "@"
$wsStocks.Cells.Item(97,10).Value = "2025-06-19"

# Row 98
$wsStocks.Cells.Item(98,10).NumberFormat = "@"
$wsStocks.Cells.Item(98,10).Value = "2025-06-19"
$wsStocks.Cells.Item(98,11).Value = 280.5

# Row 101
$wsStocks.Cells.Item(101,9).Value = 200
$wsStocks.Cells.Item(101,10).NumberFormat = "@"
$wsStocks.Cells.Item(101,10).Value = "2025-06-19"
$wsStocks.Cells.Item(101,11).Value = 240

# Row 103
$wsStocks.Cells.Item(103,8).Value = $true
$wsStocks.Cells.Item(103,9).Value = 539
$wsStocks.Cells.Item(103,10).NumberFormat = "@"
$wsStocks.Cells.Item(103,10).Value = "2025-07-25"
$wsStocks.Cells.Item(103,11).Value = 450.84

# Row 107
$wsStocks.Cells.Item(107,10).NumberFormat = "@"
$wsStocks.Cells.Item(107,10).Value = "2025-06-19"

# Row 110
$wsStocks.Cells.Item(110,7).Value = $true
$wsStocks.Cells.Item(110,8).Value = $true
$wsStocks.Cells.Item(110,9).Value = 350
$wsStocks.Cells.Item(110,10).NumberFormat = "@"
$wsStocks.Cells.Item(110,10).Value = "2025-07-25T02:18:06.713Z"
$wsStocks.Cells.Item(110,11).Value = 325.55

# Row 115
$wsStocks.Cells.Item(115,10).NumberFormat = "@"
$wsStocks.Cells.Item(115,10).Value = "2025-06-19"

# Row 116
$wsStocks.Cells.Item(116,7).Value = $true
$wsStocks.Cells.Item(116,10).NumberFormat = "@"
$wsStocks.Cells.Item(116,10).Value = "2025-06-19"

# Row 119
$wsStocks.Cells.Item(119,8).Value = $false
$wsStocks.Cells.Item(119,9).ClearContents()
$wsStocks.Cells.Item(119,10).ClearContents()
$wsStocks.Cells.Item(119,11).ClearContents()

# Row 121
$wsStocks.Cells.Item(121,7).Value = $true
$wsStocks.Cells.Item(121,8).Value = $true
$wsStocks.Cells.Item(121,9).Value = 257
$wsStocks.Cells.Item(121,10).NumberFormat = "@"
$wsStocks.Cells.Item(121,10).Value = "2025-07-14"
$wsStocks.Cells.Item(121,11).Value = 218.45

# Row 122
$wsStocks.Cells.Item(122,8).Value = $false
$wsStocks.Cells.Item(122,9).ClearContents()
$wsStocks.Cells.Item(122,10).ClearContents()
$wsStocks.Cells.Item(122,11).ClearContents()

# Row 126
$wsStocks.Cells.Item(126,7).Value = $true

# Row 131
$wsStocks.Cells.Item(131,10).NumberFormat = "@"
$wsStocks.Cells.Item(131,10).Value = "2025-06-19"
$wsStocks.Cells.Item(131,11).Value = 300.75

# Row 132
$wsStocks.Cells.Item(132,9).Value = 200.78
$wsStocks.Cells.Item(132,10).NumberFormat = "@"
$wsStocks.Cells.Item(132,10).Value = "2025-06-19"
$wsStocks.Cells.Item(132,11).Value = 340.46

# Row 136
$wsStocks.Cells.Item(136,10).NumberFormat = "@"
$wsStocks.Cells.Item(136,10).Value = "2025-06-19"

# Row 137
$wsStocks.Cells.Item(137,10).NumberFormat = "@"
$wsStocks.Cells.Item(137,10).Value = "2025-06-19"

# Row 141
$wsStocks.Cells.Item(141,10).NumberFormat = "@"
$wsStocks.Cells.Item(141,10).Value = "2025-06-19"
$wsStocks.Cells.Item(141,11).Value = 190.91

# Row 148
$wsStocks.Cells.Item(148,10).NumberFormat = "@"
$wsStocks.Cells.Item(148,10).Value = "2025-06-19"

# Row 154
$wsStocks.Cells.Item(154,8).Value = $true
$wsStocks.Cells.Item(154,9).Value = 694
$wsStocks.Cells.Item(154,10).NumberFormat = "@"
$wsStocks.Cells.Item(154,10).Value = "2025-07-14"
$wsStocks.Cells.Item(154,11).Value = 537.2

# Row 155
$wsStocks.Cells.Item(155,10).NumberFormat = "@"
$wsStocks.Cells.Item(155,10).Value = "2025-06-19"

# Row 156
$wsStocks.Cells.Item(156,8).Value = $false
$wsStocks.Cells.Item(156,9).ClearContents()
$wsStocks.Cells.Item(156,10).ClearContents()
$wsStocks.Cells.Item(156,11).ClearContents()

# Row 157
$wsStocks.Cells.Item(157,9).Value = 490
$wsStocks.Cells.Item(157,10).NumberFormat = "@"
$wsStocks.Cells.Item(157,10).Value = "2025-07-25"
$wsStocks.Cells.Item(157,11).Value = 397.8

# Row 158
$wsStocks.Cells.Item(158,10).NumberFormat = "@"
$wsStocks.Cells.Item(158,10).Value = "2025-06-19"

# Row 161
$wsStocks.Cells.Item(161,8).Value = $false
$wsStocks.Cells.Item(161,9).ClearContents()
$wsStocks.Cells.Item(161,10).ClearContents()
$wsStocks.Cells.Item(161,11).ClearContents()

# Row 164
$wsStocks.Cells.Item(164,8).Value = $true
$wsStocks.Cells.Item(164,9).Value = 503.9
$wsStocks.Cells.Item(164,10).NumberFormat = "@"
$wsStocks.Cells.Item(164,10).Value = "2025-07-14"
$wsStocks.Cells.Item(164,11).Value = 428.31

# Row 165
$wsStocks.Cells.Item(165,8).Value = $true
$wsStocks.Cells.Item(165,9).Value = 529
$wsStocks.Cells.Item(165,10).NumberFormat = "@"
$wsStocks.Cells.Item(165,10).Value = "2025-07-14"
$wsStocks.Cells.Item(165,11).Value = 449.65

# Row 166
$wsStocks.Cells.Item(166,10).NumberFormat = "@"
$wsStocks.Cells.Item(166,10).Value = "2025-06-19"

# Row 182
$wsStocks.Cells.Item(182,7).Value = $true

# Row 187
$wsStocks.Cells.Item(187,8).Value = $false
$wsStocks.Cells.Item(187,9).ClearContents()
$wsStocks.Cells.Item(187,10).ClearContents()
$wsStocks.Cells.Item(187,11).ClearContents()

# Row 202
$wsStocks.Cells.Item(202,8).Value = $false
$wsStocks.Cells.Item(202,9).ClearContents()
$wsStocks.Cells.Item(202,10).ClearContents()
$wsStocks.Cells.Item(202,11).ClearContents()

# Row 203
$wsStocks.Cells.Item(203,7).Value = $true

# Row 204
$wsStocks.Cells.Item(204,8).Value = $false
$wsStocks.Cells.Item(204,9).ClearContents()
$wsStocks.Cells.Item(204,10).ClearContents()
$wsStocks.Cells.Item(204,11).ClearContents()

# Row 220
$wsStocks.Cells.Item(220,8).Value = $false
$wsStocks.Cells.Item(220,9).ClearContents()
$wsStocks.Cells.Item(220,10).ClearContents()
$wsStocks.Cells.Item(220,11).ClearContents()

# Row 226
$wsStocks.Cells.Item(226,10).NumberFormat = "@"
$wsStocks.Cells.Item(226,10).Value = "2025-06-19"

# Row 233
$wsStocks.Cells.Item(233,10).NumberFormat = "@"
$wsStocks.Cells.Item(233,10).Value = "2025-06-19"

# Row 244
$wsStocks.Cells.Item(244,12).ClearContents()

# Row 245
$wsStocks.Cells.Item(245,12).ClearContents()

# Row 246
$wsStocks.Cells.Item(246,12).ClearContents()

# Row 247
$wsStocks.Cells.Item(247,12).ClearContents()

# Row 248
$wsStocks.Cells.Item(248,1).Value = "TTL"
$wsStocks.Cells.Item(248,2).Value = 400
$wsStocks.Cells.Item(248,3).Value = 500
$wsStocks.Cells.Item(248,4).Value = 600
$wsStocks.Cells.Item(248,5).Value = 700
$wsStocks.Cells.Item(248,6).Value = "OTHERS"
$wsStocks.Cells.Item(248,7).Value = $false
$wsStocks.Cells.Item(248,8).Value = $false
$wsStocks.Cells.Item(248,12).NumberFormat = "@"
$wsStocks.Cells.Item(248,12).Value = "2025-07-20T10:25:13.270Z"

# Row 249
$wsStocks.Cells.Item(249,1).Value = "SANVI"
$wsStocks.Cells.Item(249,2).Value = 300
$wsStocks.Cells.Item(249,3).Value = 400
$wsStocks.Cells.Item(249,4).Value = 500
$wsStocks.Cells.Item(249,5).Value = 900
$wsStocks.Cells.Item(249,6).Value = "HYDROPOWER"
$wsStocks.Cells.Item(249,7).Value = $false
$wsStocks.Cells.Item(249,8).Value = $false
$wsStocks.Cells.Item(249,12).NumberFormat = "@"
$wsStocks.Cells.Item(249,12).Value = "2025-07-20T10:25:38.664Z"

#### SupportHistory sheet (sheet2) changes ####
# Row 63
$wsSupport.Cells.Item(63,1).Value = "GLBSL"
$wsSupport.Cells.Item(63,2).Value = "support2"
$wsSupport.Cells.Item(63,3).NumberFormat = "@"
$wsSupport.Cells.Item(63,3).Value = "2025-07-02T13:36:03.613Z"

# Row 64
$wsSupport.Cells.Item(64,1).Value = "MLBS"
$wsSupport.Cells.Item(64,2).Value = "support1"
$wsSupport.Cells.Item(64,3).NumberFormat = "@"
$wsSupport.Cells.Item(64,3).Value = "2025-07-02T13:36:03.614Z"

# Row 65
$wsSupport.Cells.Item(65,1).Value = "SHLB"
$wsSupport.Cells.Item(65,2).Value = "support1"
$wsSupport.Cells.Item(65,3).NumberFormat = "@"
$wsSupport.Cells.Item(65,3).Value = "2025-07-02T13:36:03.614Z"

# Row 66
$wsSupport.Cells.Item(66,1).Value = "NMFBS"
$wsSupport.Cells.Item(66,2).Value = "support1"
$wsSupport.Cells.Item(66,3).NumberFormat = "@"
$wsSupport.Cells.Item(66,3).Value = "2025-07-02T13:36:03.615Z"

# Row 67
$wsSupport.Cells.Item(67,1).Value = "EHPL"
$wsSupport.Cells.Item(67,2).Value = "support1"
$wsSupport.Cells.Item(67,3).NumberFormat = "@"
$wsSupport.Cells.Item(67,3).Value = "2025-07-02T13:36:03.616Z"

# Row 68
$wsSupport.Cells.Item(68,1).Value = "MANDU"
$wsSupport.Cells.Item(68,2).Value = "support1"
$wsSupport.Cells.Item(68,3).NumberFormat = "@"
$wsSupport.Cells.Item(68,3).Value = "2025-07-02T13:36:03.616Z"

# Row 69
$wsSupport.Cells.Item(69,1).Value = "MKJC"
$wsSupport.Cells.Item(69,2).Value = "support1"
$wsSupport.Cells.Item(69,3).NumberFormat = "@"
$wsSupport.Cells.Item(69,3).Value = "2025-07-02T13:36:03.617Z"

# Row 70
$wsSupport.Cells.Item(70,1).Value = "NHPC"
$wsSupport.Cells.Item(70,2).Value = "support1"
$wsSupport.Cells.Item(70,3).NumberFormat = "@"
$wsSupport.Cells.Item(70,3).Value = "2025-07-02T13:36:03.617Z"

# Row 71
$wsSupport.Cells.Item(71,1).Value = "NLG"
$wsSupport.Cells.Item(71,2).Value = "support2"
$wsSupport.Cells.Item(71,3).NumberFormat = "@"
$wsSupport.Cells.Item(71,3).Value = "2025-07-02T13:36:03.618Z"

# Row 72
$wsSupport.Cells.Item(72,1).Value = "CFCL"
$wsSupport.Cells.Item(72,2).Value = "support1"
$wsSupport.Cells.Item(72,3).NumberFormat = "@"
$wsSupport.Cells.Item(72,3).Value = "2025-07-02T13:36:03.619Z"

# Row 73
$wsSupport.Cells.Item(73,1).Value = "SHIVM"
$wsSupport.Cells.Item(73,2).Value = "support1"
$wsSupport.Cells.Item(73,3).NumberFormat = "@"
$wsSupport.Cells.Item(73,3).Value = "2025-07-02T13:36:03.619Z"

# Row 74
$wsSupport.Cells.Item(74,1).Value = "GBBL"
$wsSupport.Cells.Item(74,2).Value = "support1"
$wsSupport.Cells.Item(74,3).NumberFormat = "@"
$wsSupport.Cells.Item(74,3).Value = "2025-07-02T13:36:03.620Z"

# Row 75
$wsSupport.Cells.Item(75,1).Value = "SNLI"
$wsSupport.Cells.Item(75,2).Value = "support2"
$wsSupport.Cells.Item(75,3).NumberFormat = "@"
$wsSupport.Cells.Item(75,3).Value = "2025-07-02T13:36:03.621Z"

# Row 76
$wsSupport.Cells.Item(76,1).Value = "NADEP"
$wsSupport.Cells.Item(76,2).Value = "support1"
$wsSupport.Cells.Item(76,3).NumberFormat = "@"
$wsSupport.Cells.Item(76,3).Value = "2025-07-20T10:19:35.208Z"

# Row 77
$wsSupport.Cells.Item(77,1).Value = "HHL"
$wsSupport.Cells.Item(77,2).Value = "support1"
$wsSupport.Cells.Item(77,3).NumberFormat = "@"
$wsSupport.Cells.Item(77,3).Value = "2025-07-20T10:19:35.209Z"

# Row 78
$wsSupport.Cells.Item(78,1).Value = "NLG"
$wsSupport.Cells.Item(78,2).Value = "support1"
$wsSupport.Cells.Item(78,3).NumberFormat = "@"
$wsSupport.Cells.Item(78,3).Value = "2025-07-20T10:19:35.209Z"

# Row 79
$wsSupport.Cells.Item(79,1).Value = "JFL"
$wsSupport.Cells.Item(79,2).Value = "support1"
$wsSupport.Cells.Item(79,3).NumberFormat = "@"
$wsSupport.Cells.Item(79,3).Value = "2025-07-20T10:19:35.209Z"

# Row 80
$wsSupport.Cells.Item(80,1).Value = "HBL"
$wsSupport.Cells.Item(80,2).Value = "support1"
$wsSupport.Cells.Item(80,3).NumberFormat = "@"
$wsSupport.Cells.Item(80,3).Value = "2025-07-20T10:19:35.209Z"

# Row 81
$wsSupport.Cells.Item(81,1).Value = "TRH"
$wsSupport.Cells.Item(81,2).Value = "support1"
$wsSupport.Cells.Item(81,3).NumberFormat = "@"
$wsSupport.Cells.Item(81,3).Value = "2025-07-25T02:17:57.727Z"

# Row 82
$wsSupport.Cells.Item(82,1).Value = "SRLI"
$wsSupport.Cells.Item(82,2).Value = "support1"
$wsSupport.Cells.Item(82,3).NumberFormat = "@"
$wsSupport.Cells.Item(82,3).Value = "2025-07-25T02:17:57.728Z"

# Row 83
$wsSupport.Cells.Item(83,1).Value = "OMPL"
$wsSupport.Cells.Item(83,2).Value = "support1"
$wsSupport.Cells.Item(83,3).NumberFormat = "@"
$wsSupport.Cells.Item(83,3).Value = "2025-07-29T15:23:55.610Z"
